# Refresh the cryptos list: update Price (D) / Volume(1h) (E) for every coin,
# and fix the PancakeSwap/Litecoin row order (rows 24-25 swapped places).
#
# Price values are strings like "69.586.63" that must stay text (Excel would
# otherwise coerce plain numeric-looking ones, e.g. "1.00" -> 1). We prefix
# them with a leading apostrophe to force text entry, then reset the cell's
# Style back to "Normal" so no extra number-format/style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''69.586.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.27%  '

# Row 3
$ws.Range("D3").Value = '''3.500.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.36%  '

# Row 4
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").Value = '''602.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.37%  '

# Row 6
$ws.Range("D6").Value = '''194.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.30%  '

# Row 7
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("D8").Value = '''0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = '''0.201'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.60%  '

# Row 10
$ws.Range("D10").Value = '''0.647'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '

# Row 11
$ws.Range("D11").Value = '''53.17'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '

# Row 12
$ws.Range("D12").Value = '''0.0000300'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.80%  '

# Row 13
$ws.Range("D13").Value = '''9.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '

# Row 14
$ws.Range("D14").Value = '''4.057.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '

# Row 15
$ws.Range("D15").Value = '''592.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.55%  '

# Row 16
$ws.Range("D16").Value = '''69.737.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '

# Row 17
$ws.Range("D17").Value = '''19.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.85%  '

# Row 18
$ws.Range("D18").Value = '''12.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.42%  '

# Row 19
$ws.Range("E19").Value = '  +2.34%  '

# Row 20
$ws.Range("D20").Value = '''3.507.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.48%  '

# Row 21
$ws.Range("D21").Value = '''0.984'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.57%  '

# Row 22
$ws.Range("D22").Value = '''18.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.38%  '

# Row 23
$ws.Range("D23").Value = '''5.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.55%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''101.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.23%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''4.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.37%  '

# Row 26
$ws.Range("E26").Value = '  +4.66%  '

# Row 27
$ws.Range("D27").Value = '''10.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.86%  '

# Row 28
$ws.Range("D28").Value = '''9.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.21%  '

# Row 29
$ws.Range("D29").Value = '''33.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.34%  '

# Row 30
$ws.Range("D30").Value = '''4.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.44%  '

# Row 31
$ws.Range("D31").Value = '''7.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.12%  '

# Row 32
$ws.Range("D32").Value = '''12.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.78%  '

# Row 33
$ws.Range("E33").Value = '  +0.48%  '

# Row 34
$ws.Range("D34").Value = '''63.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("D35").Value = '''3.733.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.30%  '

# Row 36
$ws.Range("D36").Value = '''3.11'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.69%  '

# Row 37
$ws.Range("D37").Value = '''0.0₃0814'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.25%  '

# Row 38
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '

# Row 39
$ws.Range("D39").Value = '''3.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.15%  '

# Row 40
$ws.Range("D40").Value = '''0.390'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.40%  '

# Row 41
$ws.Range("D41").Value = '''36.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.66%  '

# Row 42
$ws.Range("D42").Value = '''485.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.12%  '

# Row 43
$ws.Range("E43").Value = '  -1.49%  '

# Row 44
$ws.Range("D44").Value = '''0.0452'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.52%  '

# Row 45
$ws.Range("D45").Value = '''0.140'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.79%  '

# Row 46
$ws.Range("D46").Value = '''2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.85%  '

# Row 47
$ws.Range("E47").Value = '  -1.27%  '

# Row 48
$ws.Range("D48").Value = '''1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '

# Row 49
$ws.Range("D49").Value = '''8.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.65%  '

# Row 50
$ws.Range("E50").Value = '  +2.79%  '

# Row 51
$ws.Range("E51").Value = '  +10.26%  '

